$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024831875083829
$ws.Range("D2").Value = 1.033231186621255
$ws.Range("E2").Value = 1.025262092286521
$ws.Range("F2").Value = 1.041822726686595
$ws.Range("I2").Value = 1.030703326759
$ws.Range("J2").Value = 1.030004285704727
$ws.Range("K2").Value = 1.036034131710098
$ws.Range("L2").Value = 1.028088178694022
$ws.Range("M2").Value = 1.044601153539867
$ws.Range("N2").Value = 1.031467009933577

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025784332309181
$ws.Range("D3").Value = 1.033955741450153
$ws.Range("E3").Value = 1.026070820503341
$ws.Range("F3").Value = 1.042782500365465
$ws.Range("I3").Value = 1.03083766256289
$ws.Range("J3").Value = 1.030595802537291
$ws.Range("K3").Value = 1.03656784919079
$ws.Range("L3").Value = 1.028704168659835
$ws.Range("M3").Value = 1.045371242308983
$ws.Range("N3").Value = 1.032059366787885

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026401156995317
$ws.Range("D4").Value = 1.034424823200633
$ws.Range("E4").Value = 1.026594950050615
$ws.Range("F4").Value = 1.043404305906185
$ws.Range("I4").Value = 1.030923291962704
$ws.Range("J4").Value = 1.030978482296347
$ws.Range("K4").Value = 1.036912777896708
$ws.Range("L4").Value = 1.029102934168904
$ws.Range("M4").Value = 1.045869688963475
$ws.Range("N4").Value = 1.032442589996093

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.026660593926928
$ws.Range("D5").Value = 1.034622082721213
$ws.Range("E5").Value = 1.026815491061498
$ws.Range("F5").Value = 1.043665895374549
$ws.Range("I5").Value = 1.030958980031071
$ws.Range("J5").Value = 1.03113934285611
$ws.Range("K5").Value = 1.03705768379023
$ws.Range("L5").Value = 1.029270616955149
$ws.Range("M5").Value = 1.046079270051868
$ws.Range("N5").Value = 1.032603678996301

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.026704161763873
$ws.Range("D6").Value = 1.034655206796365
$ws.Range("E6").Value = 1.026852532377589
$ws.Range("F6").Value = 1.043709828054054
$ws.Range("I6").Value = 1.030964953989339
$ws.Range("J6").Value = 1.031166350967594
$ws.Range("K6").Value = 1.037082008107328
$ws.Range("L6").Value = 1.029298774051724
$ws.Range("M6").Value = 1.046114461602977
$ws.Range("N6").Value = 1.032630725462402

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026404623118877
$ws.Range("D7").Value = 1.034427458767643
$ws.Range("E7").Value = 1.026597896158217
$ws.Range("F7").Value = 1.043407800562427
$ws.Range("I7").Value = 1.030923770049272
$ws.Range("J7").Value = 1.030980631795006
$ws.Range("K7").Value = 1.036914714539
$ws.Range("L7").Value = 1.029105174591695
$ws.Range("M7").Value = 1.045872489263758
$ws.Range("N7").Value = 1.032444742547286

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025153653986139
$ws.Range("D8").Value = 1.033476001461907
$ws.Range("E8").Value = 1.025535233454138
$ws.Range("F8").Value = 1.042146926880958
$ws.Range("I8").Value = 1.030748993953574
$ws.Range("J8").Value = 1.030204205744204
$ws.Range("K8").Value = 1.036214590811526
$ws.Range("L8").Value = 1.028296317736652
$ws.Range("M8").Value = 1.0448613773044
$ws.Range("N8").Value = 1.031667213882435

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022953319441341
$ws.Range("D9").Value = 1.031801361269094
$ws.Range("E9").Value = 1.023669089500834
$ws.Range("F9").Value = 1.039931048184654
$ws.Range("I9").Value = 1.030431123092232
$ws.Range("J9").Value = 1.028835541406988
$ws.Range("K9").Value = 1.034977696654702
$ws.Range("L9").Value = 1.026872424053712
$ws.Range("M9").Value = 1.043080855478284
$ws.Range("N9").Value = 1.030296605884872

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021489195659868
$ws.Range("D10").Value = 1.030686328756293
$ws.Range("E10").Value = 1.022429378653272
$ws.Range("F10").Value = 1.038457872607516
$ws.Range("I10").Value = 1.030212589678071
$ws.Range("J10").Value = 1.027922817678369
$ws.Range("K10").Value = 1.034151020502537
$ws.Range("L10").Value = 1.025924174623114
$ws.Range("M10").Value = 1.041894712416805
$ws.Range("N10").Value = 1.029382585983867

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020855880978255
$ws.Range("D11").Value = 1.030203855558442
$ws.Range("E11").Value = 1.02189362689345
$ws.Range("F11").Value = 1.037820955284493
$ws.Range("I11").Value = 1.030116398891163
$ws.Range("J11").Value = 1.027527543406635
$ws.Range("K11").Value = 1.033792580228143
$ws.Range("L11").Value = 1.025513825171196
$ws.Range("M11").Value = 1.041381321356535
$ws.Range("N11").Value = 1.028986750377329

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020620739786641
$ws.Range("D12").Value = 1.030024696485486
$ws.Range("E12").Value = 1.021694784091993
$ws.Range("F12").Value = 1.03758452401427
$ws.Range("I12").Value = 1.030080434828302
$ws.Range("J12").Value = 1.02738071301633
$ws.Range("K12").Value = 1.033659367843922
$ws.Range("L12").Value = 1.025361441496992
$ws.Range("M12").Value = 1.041190658813434
$ws.Range("N12").Value = 1.028839711471028

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020671173828179
$ws.Range("D13").Value = 1.030063124296227
$ws.Range("E13").Value = 1.021737429319965
$ws.Range("F13").Value = 1.037635232597297
$ws.Range("I13").Value = 1.030088159846301
$ws.Range("J13").Value = 1.02741220897739
$ws.Range("K13").Value = 1.03368794558111
$ws.Range("L13").Value = 1.025394126559377
$ws.Range("M13").Value = 1.041231555029416
$ws.Range("N13").Value = 1.028871252159965

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020836442087564
$ws.Range("D14").Value = 1.030189045125302
$ws.Range("E14").Value = 1.021877187229019
$ws.Range("F14").Value = 1.037801408771338
$ws.Range("I14").Value = 1.030113430871737
$ws.Range("J14").Value = 1.027515406517142
$ws.Range("K14").Value = 1.033781570307829
$ws.Range("L14").Value = 1.02550122830035
$ws.Range("M14").Value = 1.041365560435719
$ws.Range("N14").Value = 1.028974596252061

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020938282603289
$ws.Range("D15").Value = 1.030266636162513
$ws.Range("E15").Value = 1.021963317831978
$ws.Range("F15").Value = 1.037903815072244
$ws.Range("I15").Value = 1.030128970123375
$ws.Range("J15").Value = 1.027578988904543
$ws.Range("K15").Value = 1.033839246120511
$ws.Range("L15").Value = 1.025567222333883
$ws.Range("M15").Value = 1.041448130098023
$ws.Range("N15").Value = 1.029038268933745

$ws.Range("B16").Value = 1.019999999999999
$ws.Range("C16").Value = 1.021531240658292
$ws.Range("D16").Value = 1.030718356263516
$ws.Range("E16").Value = 1.022464956993993
$ws.Range("F16").Value = 1.038500163419746
$ws.Range("I16").Value = 1.03021894062754
$ws.Range("J16").Value = 1.027949049567987
$ws.Range("K16").Value = 1.034174798891456
$ws.Range("L16").Value = 1.025951413492449
$ws.Range("M16").Value = 1.041928789158675
$ws.Range("N16").Value = 1.029408855125777

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021903365023833
$ws.Range("D17").Value = 1.031001801222858
$ws.Range("E17").Value = 1.02277990436205
$ws.Range("F17").Value = 1.038874499609033
$ws.Range("I17").Value = 1.030274958278233
$ws.Range("J17").Value = 1.028181163731453
$ws.Range("K17").Value = 1.034385153403671
$ws.Range("L17").Value = 1.026192473776715
$ws.Range("M17").Value = 1.042230352753197
$ws.Range("N17").Value = 1.029641298917979

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022120482550271
$ws.Range("D18").Value = 1.031167162896679
$ws.Range("E18").Value = 1.022963709171941
$ws.Range("F18").Value = 1.03909293779249
$ws.Range("I18").Value = 1.030307481452356
$ws.Range("J18").Value = 1.028316546253746
$ws.Range("K18").Value = 1.034507802860246
$ws.Range("L18").Value = 1.026333104101327
$ws.Range("M18").Value = 1.042406270525122
$ws.Range("N18").Value = 1.029776873698983

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022194524791844
$ws.Range("D19").Value = 1.031223552513443
$ws.Range("E19").Value = 1.023026399010706
$ws.Range("F19").Value = 1.039167435508425
$ws.Range("I19").Value = 1.03031854539104
$ws.Range("J19").Value = 1.028362707193333
$ws.Range("K19").Value = 1.034549615191248
$ws.Range("L19").Value = 1.026381059437829
$ws.Range("M19").Value = 1.042466257435847
$ws.Range("N19").Value = 1.029823100192399

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021863433008379
$ws.Range("D20").Value = 1.030971386835295
$ws.Range("E20").Value = 1.022746102998591
$ws.Range("F20").Value = 1.03883432712357
$ws.Range("I20").Value = 1.03026896372671
$ws.Range("J20").Value = 1.028156260665702
$ws.Range("K20").Value = 1.034362589185725
$ws.Range("L20").Value = 1.02616660780618
$ws.Range("M20").Value = 1.042197995670332
$ws.Range("N20").Value = 1.029616360487018

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020787771934961
$ws.Range("D21").Value = 1.030151963106685
$ws.Range("E21").Value = 1.021836027601294
$ws.Range("F21").Value = 1.037752469928622
$ws.Range("I21").Value = 1.030105995654471
$ws.Range("J21").Value = 1.027485017625231
$ws.Range("K21").Value = 1.033754002130644
$ws.Range("L21").Value = 1.025469688439461
$ws.Range("M21").Value = 1.041326098242125
$ws.Range("N21").Value = 1.028944164204438

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020112040571534
$ws.Range("D22").Value = 1.029637065952014
$ws.Range("E22").Value = 1.021264749141932
$ws.Range("F22").Value = 1.037073121370967
$ws.Range("I22").Value = 1.03000217454132
$ws.Range("J22").Value = 1.027062935081365
$ws.Range("K22").Value = 1.033370944924226
$ws.Range("L22").Value = 1.025031730141992
$ws.Range("M22").Value = 1.040778097615468
$ws.Range("N22").Value = 1.028521482254946

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020470203347482
$ws.Range("D23").Value = 1.029909993106628
$ws.Range("E23").Value = 1.021567506819513
$ws.Range("F23").Value = 1.037433175120192
$ws.Range("I23").Value = 1.030057340499509
$ws.Range("J23").Value = 1.027286693020445
$ws.Range("K23").Value = 1.033574049735322
$ws.Range("L23").Value = 1.025263878690141
$ws.Range("M23").Value = 1.041068584147891
$ws.Range("N23").Value = 1.028745557955965

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021881476370583
$ws.Range("D24").Value = 1.030985129685963
$ws.Range("E24").Value = 1.022761376066791
$ws.Range("F24").Value = 1.038852479049545
$ws.Range("I24").Value = 1.03027167287304
$ws.Range("J24").Value = 1.028167513306993
$ws.Range("K24").Value = 1.034372785129155
$ws.Range("L24").Value = 1.026178295451
$ws.Range("M24").Value = 1.042212616379028
$ws.Range("N24").Value = 1.029627629108351

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023521675153259
$ws.Range("D25").Value = 1.032234056278317
$ws.Range("E25").Value = 1.024150765375778
$ws.Range("F25").Value = 1.04050319336801
$ws.Range("I25").Value = 1.030514469410416
$ws.Range("J25").Value = 1.029189427618647
$ws.Range("K25").Value = 1.035297834642232
$ws.Range("L25").Value = 1.027240360642637
$ws.Range("M25").Value = 1.043541014840851
$ws.Range("N25").Value = 1.030650994655543
